# "Switched to using Croutons instead of Toasts."
#
# Fills in the previously-blank time-log entry on row 82 of Sheet1
# (Sheet2's category roll-up and the pie chart derive from this data via
# formulas, so they pick the new entry up automatically on recalculation).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A82").Value = 41924
$ws.Range("B82").Value = 0.7006944444444444
$ws.Range("C82").Value = 0.8208333333333333
$ws.Range("D82").Value = 60
$ws.Range("E82").Formula = "=IF(AND(NOT(ISBLANK(B82)),NOT(ISBLANK(C82))), (C82-B82) * 24 - D82/60, """")"
$ws.Range("F82").Value = "Coding"

# Matches the author's recorded cursor position after making the edit.
$null = $ws.Range("D83").Select()

$excel.CalculateFullRebuild()
